$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: recalculated date + new comment "Counts below detection"
$ws.Range("L17").Value = 42507
$ws.Range("N17").Value = "Counts below detection"

# Row 18: recalculated date
$ws.Range("L18").Value = 42507

# Row 19: recalculated Cs / sCs
$ws.Range("E19").Value = 1296.4122170718033
$ws.Range("F19").Value = 14.643579148855034

# Row 20: recalculated Cs, sCs, fSorb, sfSorb, date
$ws.Range("E20").Value = 2623.7646921398959
$ws.Range("F20").Value = 0
$ws.Range("H20").Value = 0.78886057483628635
$ws.Range("I20").Value = 0
$ws.Range("L20").Value = 42507

# Row 21: recalculated Cw, sCw, Cs, sCs, fSorb, sfSorb, pH, spH, date
$ws.Range("C21").Value = 0.82945425208510437
$ws.Range("D21").Value = 0.51618184155151592
$ws.Range("E21").Value = 13907.345389394934
$ws.Range("F21").Value = 1720.6061385050698
$ws.Range("H21").Value = 0.83416413640724907
$ws.Range("I21").Value = 0.10320214918352018
$ws.Range("J21").Value = 8.9350000000000005
$ws.Range("K21").Value = 0.04725815626252589
$ws.Range("L21").Value = 42507

# Update the active selection to match the reviewed range (C21:D21)
$ws.Activate() | Out-Null
$ws.Range("C21:D21").Select() | Out-Null
